$d = $word.ActiveDocument

# Locate the paragraph that contains "Tests de la version courante"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Tests de la version courante") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Tests de la version courante'"
}

# Insert a new paragraph right after it, inheriting its paragraph formatting
$newRange = $target.Range.InsertParagraphAfter()

# The newly created paragraph is the one following $target
$newPara = $target.Next()
$newPara.Range.Text = "Identifier les TEC facturés"
